$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 0.0001388888888888889
$ws.Range("K2").Value = 966
$ws.Range("L2").Value = 0.001932
